$d = $word.ActiveDocument

# Locate the "Restrições" / "Razão (lógica)" table; the paragraph that
# immediately follows it is where the replacement text belongs once the
# table itself is gone.
$tbl = $d.Tables.Item(1)
$insertPos = $tbl.Range.End

# Insert the new run — matching the document's usual character formatting
# (Arial, black/text1 colour, 12pt) — right after the table, before the
# table is removed (ranges computed from the table become unreliable once
# the table itself has been deleted).
$target = $d.Range($insertPos, $insertPos)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Não há restrições.</w:t></w:r></w:p>'
$target.InsertXML($xml)

# Now remove the whole table.
$d.Tables.Item(1).Delete()
